$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 9722.75
$ws.Range("I20").Value = 8330.333000000001
$ws.Range("K20").Value = 8330.333000000001
$ws.Range("M20").Value = -8100.333000000001
$ws.Range("H35").Value = 9722.75
$ws.Range("I35").Value = 8330.333000000001
$ws.Range("K35").Value = 8330.333000000001
$ws.Range("M35").Value = -7951.333000000001
$ws.Range("H69").Value = 3497.25
$ws.Range("J69").Value = 3000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10748
$ws.Range("H72").Value = 3497.25
$ws.Range("J72").Value = 3000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35736
$ws.Range("H107").Value = 1698.25
$ws.Range("I107").Value = 1599.5
$ws.Range("J107").Value = 1797
$ws.Range("K107").Value = 1599.5
$ws.Range("L107").Value = 1797
$ws.Range("M107").Value = 320.5
$ws.Range("N107").Value = -5637
$ws.Range("H112").Value = 1318.4
$ws.Range("J112").Value = 1318.4
$ws.Range("L112").Value = 3955.2
$ws.Range("N112").Value = -6171.200000000001
$ws.Range("H125").Value = 880
$ws.Range("J125").Value = 880
$ws.Range("L125").Value = 7920
$ws.Range("N125").Value = -12840
$ws.Range("H132").Value = 37903
$ws.Range("I132").Value = 39195.926
$ws.Range("K132").Value = 117587.778
$ws.Range("M132").Value = -115057.778
$ws.Range("H137").Value = 1930
$ws.Range("I137").Value = 1930
$ws.Range("K137").Value = 5790
$ws.Range("M137").Value = -3240
$ws.Range("H141").Value = 4254.9287
$ws.Range("I141").Value = 3767.1
$ws.Range("K141").Value = 11301.3
$ws.Range("M141").Value = -6121.299999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7880
$ws.Range("I2").Value = 8638
$ws.Range("K2").Value = 8638
$ws.Range("M2").Value = -8525
$ws.Range("H101").Value = 44900
$ws.Range("J101").Value = 44900
$ws.Range("L101").Value = 44900
$ws.Range("N101").Value = -51390
$ws.Range("H110").Value = 5060.8887
$ws.Range("I110").Value = 4576.643
$ws.Range("J110").Value = 6755.75
$ws.Range("K110").Value = 4576.643
$ws.Range("L110").Value = 6755.75
$ws.Range("M110").Value = -2531.643
$ws.Range("N110").Value = -10845.75
$ws.Range("H116").Value = 7880
$ws.Range("I116").Value = 8638
$ws.Range("K116").Value = 8638
$ws.Range("M116").Value = -6344

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7880
$ws.Range("I3").Value = 8638
$ws.Range("K3").Value = 8638
$ws.Range("M3").Value = -8524
$ws.Range("H105").Value = 3003.1667
$ws.Range("I105").Value = 2933.5
$ws.Range("J105").Value = 3142.5
$ws.Range("K105").Value = 2933.5
$ws.Range("L105").Value = 3142.5
$ws.Range("M105").Value = -1186.5
$ws.Range("N105").Value = -6636.5
$ws.Range("H123").Value = 12225.333
$ws.Range("J123").Value = 12225.333
$ws.Range("L123").Value = 12225.333
$ws.Range("N123").Value = -22025.333
$ws.Range("H134").Value = 14708427
$ws.Range("I134").Value = 16669177
$ws.Range("J134").Value = 2803
$ws.Range("K134").Value = 50007531
$ws.Range("L134").Value = 8409
$ws.Range("M134").Value = -50004996
$ws.Range("N134").Value = -13479

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2157
$ws.Range("I31").Value = 1250.9375
$ws.Range("K31").Value = 1250.9375
$ws.Range("M31").Value = -955.9375
$ws.Range("H34").Value = 2157
$ws.Range("I34").Value = 1250.9375
$ws.Range("K34").Value = 1250.9375
$ws.Range("M34").Value = -1048.9375
$ws.Range("H43").Value = 17946
$ws.Range("J43").Value = 17946
$ws.Range("L43").Value = 17946
$ws.Range("N43").Value = -18314
$ws.Range("H99").Value = 2481
$ws.Range("I99").Value = 2538.25
$ws.Range("K99").Value = 2538.25
$ws.Range("M99").Value = -1040.25
$ws.Range("H101").Value = 17946
$ws.Range("J101").Value = 17946
$ws.Range("L101").Value = 17946
$ws.Range("N101").Value = -24436
$ws.Range("H105").Value = 734.8333
$ws.Range("I105").Value = 529.6667
$ws.Range("K105").Value = 529.6667
$ws.Range("M105").Value = 1217.3333
$ws.Range("H122").Value = 27238.2
$ws.Range("I122").Value = 2798.6
$ws.Range("J122").Value = 51677.8
$ws.Range("K122").Value = 8395.799999999999
$ws.Range("L122").Value = 155033.4
$ws.Range("M122").Value = -5945.799999999999
$ws.Range("N122").Value = -159933.4
$ws.Range("H126").Value = 2481
$ws.Range("I126").Value = 2538.25
$ws.Range("K126").Value = 7614.75
$ws.Range("M126").Value = -5144.75
$ws.Range("H134").Value = 1993.92
$ws.Range("I134").Value = 1977.0625
$ws.Range("K134").Value = 5931.1875
$ws.Range("M134").Value = -3396.1875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 770.1111
$ws.Range("I7").Value = 156.8
$ws.Range("J7").Value = 1536.75
$ws.Range("K7").Value = 470.4
$ws.Range("L7").Value = 4610.25
$ws.Range("M7").Value = -358.4
$ws.Range("N7").Value = -4834.25
$ws.Range("H35").Value = 5719
$ws.Range("I35").Value = 2200
$ws.Range("J35").Value = 10997.5
$ws.Range("K35").Value = 6600
$ws.Range("L35").Value = 32992.5
$ws.Range("M35").Value = -6312
$ws.Range("N35").Value = -33568.5
$ws.Range("H39").Value = 8672.909
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 8672.909
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 26018.727
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -26606.727
$ws.Range("H76").Value = 7498.8
$ws.Range("I76").Value = 3747
$ws.Range("K76").Value = 11241
$ws.Range("M76").Value = -10858
$ws.Range("H79").Value = 7498.8
$ws.Range("I79").Value = 3747
$ws.Range("K79").Value = 11241
$ws.Range("M79").Value = -9915
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 4622.625
$ws.Range("K80").Value = 13500
$ws.Range("L80").Value = 13867.875
$ws.Range("M80").Value = -12564
$ws.Range("N80").Value = -15739.875
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 4622.625
$ws.Range("K83").Value = 40500
$ws.Range("L83").Value = 41603.625
$ws.Range("M83").Value = -35820
$ws.Range("N83").Value = -50963.625
$ws.Range("H131").Value = 866941.2
$ws.Range("J131").Value = 1070776.9
$ws.Range("L131").Value = 3212330.7
$ws.Range("N131").Value = -3222410.7
$ws.Range("H133").Value = 5885.5713
$ws.Range("I133").Value = 5885.5713
$ws.Range("K133").Value = 17656.7139
$ws.Range("M133").Value = -12596.7139
$ws.Range("H134").Value = 13376.131
$ws.Range("I134").Value = 7530.2
$ws.Range("K134").Value = 22590.6
$ws.Range("M134").Value = -17520.6
$ws.Range("H139").Value = 3822.1428
$ws.Range("I139").Value = 3065
$ws.Range("K139").Value = 9195
$ws.Range("M139").Value = -4055
$ws.Range("H140").Value = 6999.4287
$ws.Range("I140").Value = 6999
$ws.Range("K140").Value = 20997
$ws.Range("M140").Value = -15817

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 24916.5
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H70").Value = 4417
$ws.Range("I70").Value = 4283
$ws.Range("K70").Value = 4283
$ws.Range("M70").Value = -4013
$ws.Range("H73").Value = 4417
$ws.Range("I73").Value = 4283
$ws.Range("K73").Value = 4283
$ws.Range("M73").Value = -3347

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2113
$ws.Range("I16").Value = 1704.6666
$ws.Range("J16").Value = 2929.6667
$ws.Range("K16").Value = 1704.6666
$ws.Range("L16").Value = 2929.6667
$ws.Range("M16").Value = -1534.6666
$ws.Range("N16").Value = -3269.6667
$ws.Range("H22").Value = 3701
$ws.Range("I22").Value = 1306.5714
$ws.Range("K22").Value = 1306.5714
$ws.Range("M22").Value = -1011.5714
$ws.Range("H27").Value = 3701
$ws.Range("I27").Value = 1306.5714
$ws.Range("K27").Value = 1306.5714
$ws.Range("M27").Value = -1199.5714
$ws.Range("H55").Value = 979.94446
$ws.Range("J55").Value = 1009.3333
$ws.Range("L55").Value = 1009.3333
$ws.Range("N55").Value = -1355.3333
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 3417.5
$ws.Range("I132").Value = 3016.4285
$ws.Range("K132").Value = 9049.2855
$ws.Range("M132").Value = -6519.2855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 19999
$ws.Range("J103").Value = 19999
$ws.Range("L103").Value = 19999
$ws.Range("N103").Value = -22343
$ws.Range("H122").Value = 2635.4285
$ws.Range("I122").Value = 2779.7
$ws.Range("K122").Value = 8339.099999999999
$ws.Range("M122").Value = -5889.099999999999
